# Gift Logs Merged - 23 April 2025
#
# Update the "NEW GIFT AMT YTD" test value on the GiftLog sheet from 1 to
# 100, and leave the active selection on the GiftLog sheet at K13 (matching
# where the author last clicked before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GiftLog")

# Column C on row 2 holds the "GiftValue" test data point (stored as text).
$ws.Range("C2").Value = "100"

# Reflect the author's final cell selection on the GiftLog sheet.
$ws.Range("K13").Select() | Out-Null
